$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The updated Price-column values below are plain numeric-looking strings
# (e.g. "235.48"). Assigning such a string straight to .Value would make
# Excel auto-convert it to a Number (exactly like typing it into a General
# cell) and so silently re-render things like "0.00000000115" as "1.15E-09".
# The source keeps these cells as literal text, so force a Text number
# format on each of them first to preserve the exact string contents.
$textPriceCells = @(
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D13",
    "D14",
    "D15",
    "D17",
    "D18",
    "D20",
    "D21",
    "D22",
    "D23",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D35",
    "D36",
    "D37",
    "D40",
    "D42",
    "D43",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (D) / Volume(1h) (E) cell text for each coin row.
$ws.Range("D2").Value = '29.259.91'
$ws.Range("E2").Value = '  -0.41%  '
$ws.Range("D3").Value = '1.830.16'
$ws.Range("E3").Value = '  -0.62%  '
$ws.Range("E4").Value = '  +0.35%  '
$ws.Range("D5").Value = '235.48'
$ws.Range("E5").Value = '  -1.54%  '
$ws.Range("D6").Value = '0.6037'
$ws.Range("E6").Value = '  -3.66%  '
$ws.Range("D7").Value = '1.003'
$ws.Range("E7").Value = '  +0.26%  '
$ws.Range("D8").Value = '0.07053'
$ws.Range("E8").Value = '  -4.96%  '
$ws.Range("D9").Value = '0.2794'
$ws.Range("E9").Value = '  -3.37%  '
$ws.Range("D10").Value = '23.57'
$ws.Range("E10").Value = '  -5.49%  '
$ws.Range("D11").Value = '0.07653'
$ws.Range("E11").Value = '  -0.85%  '
$ws.Range("D12").Value = '1.824.95'
$ws.Range("E12").Value = '  -0.79%  '
$ws.Range("D13").Value = '4.793'
$ws.Range("E13").Value = '  -3.69%  '
$ws.Range("D14").Value = '0.6290'
$ws.Range("E14").Value = '  -6.64%  '
$ws.Range("D15").Value = '0.000009886'
$ws.Range("E15").Value = '  -3.98%  '
$ws.Range("D16").Value = '2.085.27'
$ws.Range("E16").Value = '  -0.28%  '
$ws.Range("D17").Value = '79.10'
$ws.Range("E17").Value = '  -3.25%  '
$ws.Range("D18").Value = '5.848'
$ws.Range("E18").Value = '  -5.79%  '
$ws.Range("D19").Value = '29.261.46'
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("D20").Value = '224.52'
$ws.Range("E20").Value = '  -3.99%  '
$ws.Range("D21").Value = '1.003'
$ws.Range("E21").Value = '  +0.24%  '
$ws.Range("D22").Value = '11.71'
$ws.Range("E22").Value = '  -4.88%  '
$ws.Range("D23").Value = '7.007'
$ws.Range("E23").Value = '  -3.94%  '
$ws.Range("E24").Value = '  +0.35%  '
$ws.Range("D25").Value = '156.03'
$ws.Range("E25").Value = '  -1.65%  '
$ws.Range("D26").Value = '0.1306'
$ws.Range("E26").Value = '  -2.78%  '
$ws.Range("D27").Value = '7.986'
$ws.Range("E27").Value = '  -6.16%  '
$ws.Range("D28").Value = '16.60'
$ws.Range("E28").Value = '  -4.02%  '
$ws.Range("D29").Value = '1.484'
$ws.Range("E29").Value = '  +1.19%  '
$ws.Range("D30").Value = '0.06451'
$ws.Range("E30").Value = '  -12.02%  '
$ws.Range("D31").Value = '1.449'
$ws.Range("E31").Value = '  -2.10%  '
$ws.Range("D32").Value = '3.840'
$ws.Range("E32").Value = '  -4.68%  '
$ws.Range("D33").Value = '3.801'
$ws.Range("E33").Value = '  -5.99%  '
$ws.Range("E34").Value = '  -2.43%  '
$ws.Range("D35").Value = '1.733'
$ws.Range("E35").Value = '  -4.63%  '
$ws.Range("D36").Value = '0.6466'
$ws.Range("E36").Value = '  -7.17%  '
$ws.Range("D37").Value = '2.545'
$ws.Range("E37").Value = '  -1.24%  '
$ws.Range("D38").Value = '1.216.58'
$ws.Range("E38").Value = '  -1.42%  '
$ws.Range("E39").Value = '  -2.60%  '
$ws.Range("D40").Value = '0.01748'
$ws.Range("E40").Value = '  -5.18%  '
$ws.Range("E41").Value = '  -5.24%  '
$ws.Range("D42").Value = '0.9011'
$ws.Range("E42").Value = '  -5.68%  '
$ws.Range("D43").Value = '1.003'
$ws.Range("E43").Value = '  +0.25%  '
$ws.Range("D44").Value = '1.999.82'
$ws.Range("E44").Value = '  +0.20%  '
$ws.Range("D45").Value = '100.19'
$ws.Range("E45").Value = '  -0.68%  '
$ws.Range("D46").Value = '62.74'
$ws.Range("E46").Value = '  -4.01%  '
$ws.Range("D47").Value = '0.00000000115'
$ws.Range("E47").Value = '  -3.93%  '
$ws.Range("D48").Value = '8.563'
$ws.Range("E48").Value = '  -3.61%  '
$ws.Range("D49").Value = '1.581'
$ws.Range("E49").Value = '  -7.79%  '
$ws.Range("D50").Value = '0.4554'
$ws.Range("E50").Value = '  -0.45%  '
$ws.Range("D51").Value = '0.05504'
$ws.Range("E51").Value = '  -2.74%  '
